# Update the cached "datetimeFigureOut" date field text from 10/24/2016 to
# 10/25/2016 everywhere it appears: on the slide master and on every slide
# layout's "Date Placeholder" shape.

$p = $ppt.ActivePresentation

$oldDate = "10/24/2016"
$newDate = "10/25/2016"
$ppPlaceholderDate = 16

function Update-DateShape($shape) {
    if ($shape.HasTextFrame -eq $false) {
        return
    }
    $isDatePlaceholder = $false
    try {
        if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $isDatePlaceholder = $true
        }
    } catch {
        $isDatePlaceholder = $false
    }
    if (-not $isDatePlaceholder) {
        if ($shape.Name -like "Date Placeholder*") {
            $isDatePlaceholder = $true
        }
    }
    if (-not $isDatePlaceholder) {
        return
    }

    $tr = $shape.TextFrame.TextRange
    if ($tr.Text -eq $oldDate) {
        $tr.Text = $newDate
    }
}

$master = $p.SlideMaster
$updated = 0

# Slide master's own "Date Placeholder" shape.
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    if (Update-DateShape $master.Shapes.Item($i)) {
        $updated = $updated + 1
    }
}

# Every slide layout owned by the master.
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        if (Update-DateShape $layout.Shapes.Item($i)) {
            $updated = $updated + 1
        }
    }
}

Write-Host "Updated date field text on $updated shape(s) ($oldDate -> $newDate)."
